$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.818.12"
$ws.Range("E2").Value = "  -0.29%  "
$ws.Range("D3").Value = "2.527.48"
$ws.Range("E3").Value = "  +3.95%  "
$ws.Range("E4").Value = "  -0.14%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "567.85"
$ws.Range("E5").Value = "  +0.89%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "147.34"
$ws.Range("E6").Value = "  +4.74%  "
$ws.Range("E7").Value = "  -0.15%  "
$ws.Range("D9").Value = "2.526.83"
$ws.Range("E9").Value = "  +4.03%  "
$ws.Range("E10").Value = "  +0.23%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.64"
$ws.Range("E11").Value = "  -1.89%  "
$ws.Range("E12").Value = "  +0.95%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "27.70"
$ws.Range("E14").Value = "  +5.68%  "
$ws.Range("D15").Value = "2.980.86"
$ws.Range("E15").Value = "  +3.90%  "
$ws.Range("D16").Value = "62.802.76"
$ws.Range("E16").Value = "  -0.13%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000141"
$ws.Range("E17").Value = "  -0.42%  "
$ws.Range("D18").Value = "2.547.76"
$ws.Range("E18").Value = "  +4.79%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.56"
$ws.Range("E19").Value = "  +3.05%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "335.51"
$ws.Range("E20").Value = "  -1.18%  "
$ws.Range("E21").Value = "  +1.67%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.76"
$ws.Range("E22").Value = "  -0.71%  "
$ws.Range("E23").Value = "  -0.05%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "65.49"
$ws.Range("E24").Value = "  +0.01%  "
$ws.Range("E25").Value = "  -3.54%  "
$ws.Range("E26").Value = "  +2.38%  "
$ws.Range("E27").Value = "  +12.81%  "
$ws.Range("E28").Value = "  +0.02%  "
$ws.Range("E29").Value = "  +1.75%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.25"
$ws.Range("E30").Value = "  +10.93%  "
$ws.Range("D31").Value = "0.0₃0810"
$ws.Range("E31").Value = "  +1.88%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.83"
$ws.Range("E32").Value = "  +0.74%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "177.71"
$ws.Range("E33").Value = "  +2.04%  "
$ws.Range("E34").Value = "  +7.12%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "409.68"
$ws.Range("E35").Value = "  +9.34%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.398"
$ws.Range("E36").Value = "  -0.06%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "18.81"
$ws.Range("E37").Value = "  +1.05%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.39"
$ws.Range("E38").Value = "  -1.62%  "
$ws.Range("E40").Value = "  +1.52%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.999"
$ws.Range("E41").Value = "  -0.32%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "39.14"
$ws.Range("E42").Value = "  -1.91%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "152.18"
$ws.Range("E43").Value = "  +4.49%  "
$ws.Range("E44").Value = "  +1.86%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "20.69"
$ws.Range("E45").Value = "  +0.88%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.604"
$ws.Range("E46").Value = "  +1.91%  "
$ws.Range("E47").Value = "  +1.44%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0520"
$ws.Range("E48").Value = "  +0.20%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0236"
$ws.Range("E49").Value = "  +5.97%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "18.33"
$ws.Range("E50").Value = "  +2.61%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.79"
$ws.Range("E51").Value = "  +2.91%  "
